$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.115.07"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "1.825.01"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.009"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3632"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07289"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8721"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.13"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("D12").Value = "1.870.55"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07623"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.343"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.474"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.37%  "
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008618"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.010"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "27.378.94"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.213"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.12%  "
$ws.Range("D24").Value = "2.096.81"
$ws.Range("E24").Value = "  +1.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.885"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.098"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.124"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08906"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.946"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.152"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.66%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7334"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.463"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.010"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.507"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.085"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05241"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01914"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.922"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.210"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5210"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1630"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.286"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4850"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.009"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.635"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06272"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.90%  "
